$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

Write-Output ("FilterMode before: " + $ws.FilterMode)
Write-Output ("AutoFilterMode before: " + $ws.AutoFilterMode)

$ws.ShowAllData()

Write-Output ("FilterMode after ShowAllData: " + $ws.FilterMode)

# 1. Unhide all rows 5-45
$ws.Rows("5:45").Hidden = $false

# 2. Remove autofilter
if ($ws.AutoFilterMode) {
    $ws.AutoFilterMode = $false
}

Write-Output ("FilterMode final: " + $ws.FilterMode)

Write-Output "done"
